# Refresh the scraped movie data: updated genre/rating scrape results in a
# revised set of rows (some movies dropped, new ones added, table grows by
# two rows) for the "Movies Playing" listing on sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Final (post-scrape) contents for A2:D27 - Title | Rating | Runtime (min) | Release Year
$data = @(
    @('Venom','7.1','112','2018'),
    @('A Star Is Born','8.5','135','2018'),
    @('A Simple Favor','7.2','117','2018'),
    @('Night School','5.5','111','2018'),
    @('The Nun','5.7','96','2018'),
    @('Smallfoot','6.8','96','2018'),
    @('The House with a Clock in Its Walls','6.3','105','2018'),
    @('Crazy Rich Asians','7.5','120','2018'),
    @('Hell Fest','5.9','89','2018'),
    @('BlacKkKlansman','7.8','135','2018'),
    @('White Boy Rick','6.6','111','2018'),
    @('Mamma Mia! Here We Go Again','7.1','114','2018'),
    @('Boku no Hero Academia the Movie','8.5','96','2018'),
    @('Christopher Robin','7.7','104','2018'),
    @('Fahrenheit 11/9','5.5','128','2018'),
    @('The Wife','7.5','100','2018'),
    @('Sui Dhaaga: Made in India','6.9','122','2018'),
    @('The Little Stranger','6.1','111','2018'),
    @('The Children Act','6.7','105','2018'),
    @('Qismat','9.4','137','2018'),
    @('The Dawn Wall','8.1','100','2018'),
    @('Hello, Mrs. Money','5.3','113','2018'),
    @('Kusama: Infinity','6.9','76','2018'),
    @('My Generation','7.2','85','2018'),
    @('Afsar','6.7','127','2018'),
    @('The Woman Who Fell to Earth','7.8','60','2018')
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $cell = $ws.Cells.Item($r, $j + 1)
        # Force text storage (matches the scraped-data columns, which are
        # all shared strings -- including numeric-looking rating/runtime/year
        # values) instead of letting COM auto-coerce to a number, then drop
        # back to the default style so no stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $rowValues[$j]
        $cell.Style = "Normal"
    }
}
